$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("Test Name"), shifting Test Name,
# Test Description, Test Type, Cloud Resource, Category, Responsibility,
# Validation Steps and USNORTHCOM Validated one column to the right.
$null = $ws.Range("E1").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("E1").Value = "Control Link"

# Restore / adjust column widths around the newly inserted column.
$ws.Columns("D").ColumnWidth = 65.2640625
$ws.Columns("E").ColumnWidth = 71.7640625
$ws.Columns("G").ColumnWidth = 65.2640625

# The AutoFilter needs to cover the new last column (M). Toggling it off
# first because Range.AutoFilter() flips the existing filter off when one
# is already active on the sheet.
$ws.AutoFilterMode = $false
$null = $ws.Range("A1:M1").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$M`$1"
    }
}

# Match the saved selection/active cell.
$null = $ws.Range("E4").Select()
